$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.230.17'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '1.787.56'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '337.86'
$ws.Range('E5').Value = '  -1.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9991'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3943'
$ws.Range('E7').Value = '  +2.81%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3439'
$ws.Range('E8').Value = '  -3.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.88'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.194'
$ws.Range('E10').Value = '  -3.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07469'
$ws.Range('E11').Value = '  -4.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9998'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.73'
$ws.Range('E13').Value = '  -3.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.459'
$ws.Range('E14').Value = '  -2.43%  '
$ws.Range('D15').Value = '1.780.84'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.103'
$ws.Range('E16').Value = '  -1.73%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001095'
$ws.Range('E17').Value = '  -3.03%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06681'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '83.64'
$ws.Range('E19').Value = '  -3.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('E21').Value = '  -0.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.514'
$ws.Range('E22').Value = '  -0.81%  '
$ws.Range('D23').Value = '27.210.00'
$ws.Range('E23').Value = '  -1.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.38'
$ws.Range('E24').Value = '  -6.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.384'
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.475'
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.510'
$ws.Range('E27').Value = '  -6.92%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.19'
$ws.Range('E28').Value = '  -4.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '157.05'
$ws.Range('E29').Value = '  +2.02%  '
$ws.Range('D30').Value = '1.981.49'
$ws.Range('E30').Value = '  -1.92%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '134.79'
$ws.Range('E31').Value = '  -1.04%  '
$ws.Range('E32').Value = '  -3.16%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.015'
$ws.Range('E33').Value = '  -6.31%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08758'
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '13.03'
$ws.Range('E35').Value = '  -6.89%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.611'
$ws.Range('E36').Value = '  -4.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.421'
$ws.Range('E37').Value = '  -4.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02387'
$ws.Range('E38').Value = '  -1.30%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6810'
$ws.Range('E39').Value = '  -3.27%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06394'
$ws.Range('E40').Value = '  -1.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2202'
$ws.Range('E41').Value = '  -2.50%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.240'
$ws.Range('E42').Value = '  -4.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.458'
$ws.Range('E43').Value = '  -6.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.35'
$ws.Range('E44').Value = '  -3.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9989'
$ws.Range('E45').Value = '  -0.12%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6412'
$ws.Range('E46').Value = '  -3.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.855'
$ws.Range('E47').Value = '  -2.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.140'
$ws.Range('E48').Value = '  -2.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '131.10'
$ws.Range('E49').Value = '  -1.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07109'
$ws.Range('E50').Value = '  -3.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.73'
$ws.Range('E51').Value = '  -2.86%  '
